$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update growth_rate (column D) values to the recomputed figures ---
$ws.Range("D2").Value = 0.04097723841932249
$ws.Range("D3").Value = 0.04269517059668157
$ws.Range("D4").Value = 0.04429979142950934
$ws.Range("D5").Value = 0.03929683307690626
$ws.Range("D6").Value = 0.03937447585678307
$ws.Range("D7").Value = 0.04033680169231903
$ws.Range("D8").Value = 0.04098002201150364
$ws.Range("D9").Value = 0.04192620461027254
$ws.Range("D10").Value = 0.03826018243295994
$ws.Range("D11").Value = 0.02764419136302141
$ws.Range("D12").Value = 0.04198000344805354
$ws.Range("D13").Value = 0.04008159921462297
$ws.Range("D14").Value = 0.1190763531796746
$ws.Range("D15").Value = 0.1406605443223997
$ws.Range("D16").Value = 0.04027771755062181
$ws.Range("D18").Value = 0.05003201250530515
$ws.Range("D19").Value = 0.115877070428584
$ws.Range("D20").Value = 0.0419079400706974
$ws.Range("D21").Value = 0.1484240591804242
$ws.Range("D22").Value = 0.04131396477675519
$ws.Range("D23").Value = 0.04047511066172597
$ws.Range("D24").Value = 0.10264326075953
$ws.Range("D25").Value = 0.03848609087056792
$ws.Range("D26").Value = 0.03422137839769211
$ws.Range("D27").Value = 0.04242393469279643
$ws.Range("D29").Value = 0.03314068796248505
$ws.Range("D30").Value = 0.03936020290587005
$ws.Range("D31").Value = 0.04118177595382631
$ws.Range("D32").Value = 0.02997771865120507
$ws.Range("D34").Value = 0.04467832015415703
$ws.Range("D35").Value = 0.04024287109049683
$ws.Range("D36").Value = 0.03453962466477918
$ws.Range("D37").Value = 0.04264325139162303
$ws.Range("D38").Value = 0.04387622852051939
$ws.Range("D39").Value = 0.0408677775956135
$ws.Range("D40").Value = 0.038678600682245
$ws.Range("D41").Value = 0.03477824664364804
$ws.Range("D42").Value = 0.04586630549120276
$ws.Range("D43").Value = 0.04187542872783673
$ws.Range("D44").Value = 0.05350581053721865
$ws.Range("D45").Value = 0.03758872027990185
$ws.Range("D48").Value = 0.04255737497074124
$ws.Range("D49").Value = 0.04076814566061661
$ws.Range("D50").Value = 0.07645569241623025
$ws.Range("D51").Value = 0.04696294142709911
$ws.Range("D52").Value = 0.05324971740219987
$ws.Range("D53").Value = 0.04456783539998523
$ws.Range("D54").Value = 0.04766440820777128
$ws.Range("D55").Value = 0.07914255984883917
$ws.Range("D56").Value = 0.05804402477717532
$ws.Range("D57").Value = 0.03351680518648362
$ws.Range("D58").Value = 0.09525931738410133
$ws.Range("D59").Value = 0.07881262454997691
$ws.Range("D60").Value = 0.08739926581235634
$ws.Range("D61").Value = 0.04799146448841025
$ws.Range("D62").Value = 0.04207969934380732
$ws.Range("D63").Value = 0.04504611458119447
$ws.Range("D64").Value = 0.04507370790845988
$ws.Range("D65").Value = 0.092353539054776
$ws.Range("D66").Value = 0.04703823414138136
$ws.Range("D67").Value = 0.09499350671566599
$ws.Range("D69").Value = 0.04393566863842759
$ws.Range("D70").Value = 0.04509846876453674
$ws.Range("D71").Value = 0.0400916932253079
$ws.Range("D72").Value = 0.04069517914201637
$ws.Range("D73").Value = 0.04336958762334631
$ws.Range("D74").Value = 0.04446683057661154
$ws.Range("D75").Value = 0.2115323652309299
$ws.Range("D76").Value = 0.1408014693890252
$ws.Range("D77").Value = 0.04665157164482927
$ws.Range("D79").Value = 0.05173489160143454
$ws.Range("D81").Value = 0.06034268724749972
$ws.Range("D82").Value = 0.1902723181581565

# --- Fix DHFR12 fragment-position label for rows 4, 6, 8 (N_term -> C_term) ---
$ws.Range("J4").Value = "DHFR12_C_term"
$ws.Range("J6").Value = "DHFR12_C_term"
$ws.Range("J8").Value = "DHFR12_C_term"

# --- Populate DHFR3 / DHFR12 fragment-position columns for rows 74-82 ---
$ws.Range("I74").Value = "DHFR3_C_term"
$ws.Range("J74").Value = "DHFR12_N_term"
$ws.Range("I75").Value = "DHFR3_C_term"
$ws.Range("J75").Value = "DHFR12_N_term"
$ws.Range("I76").Value = "DHFR3_C_term"
$ws.Range("J76").Value = "DHFR12_N_term"
$ws.Range("I77").Value = "DHFR3_C_term"
$ws.Range("J77").Value = "DHFR12_N_term"
$ws.Range("I78").Value = "DHFR3_C_term"
$ws.Range("J78").Value = "DHFR12_N_term"
$ws.Range("I79").Value = "DHFR3_C_term"
$ws.Range("J79").Value = "DHFR12_N_term"
$ws.Range("I80").Value = "DHFR3_C_term"
$ws.Range("J80").Value = "DHFR12_N_term"
$ws.Range("I81").Value = "DHFR3_C_term"
$ws.Range("J81").Value = "DHFR12_N_term"
$ws.Range("I82").Value = "DHFR3_C_term"
$ws.Range("J82").Value = "DHFR12_N_term"

Write-Host "Applied growth_rate precision updates and DHFR fragment-position grouping."
